$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.906.86"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.637.41"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.385"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").Value = "3.110.90"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "63.797.66"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "2.633.66"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "344.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.58%  "
$ws.Range("E25").Value = "  +5.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "583.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("E29").Value = "  +5.02%  "
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.404"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "153.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "162.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.24%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0590"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.634"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0249"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.793"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.31%  "
